$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 27.8037835
$ws.Range("H2").Value = 55.607567
$ws.Range("I2").Value = 0.009277807097986399
$ws.Range("J2").Value = 0.00621894358432765
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 0.9094061507180001
$ws.Range("R2").Value = 3.637624602872001
$ws.Range("S2").Value = 0.009277807097986399
$ws.Range("T2").Value = 0.00621894358432765

# Row 3
$ws.Range("I3").Value = 0.001505808648764821
$ws.Range("J3").Value = 0.001514022484498545
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 0.1475986332253333
$ws.Range("R3").Value = 0.8855917993519999
$ws.Range("S3").Value = 0.001505808648764821
$ws.Range("T3").Value = 0.001514022484498545

# Row 4
$ws.Range("G4").Value = 897.2237039999999
$ws.Range("H4").Value = 2691.671112
$ws.Range("I4").Value = 0.2993933703106574
$ws.Range("J4").Value = 0.3010264914681929
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 29.346392910432
$ws.Range("R4").Value = 176.078357462592
$ws.Range("S4").Value = 0.2993933703106574
$ws.Range("T4").Value = 0.3010264914681929

# Row 5
$ws.Range("G5").Value = 20.9707925
$ws.Range("H5").Value = 41.941585
$ws.Range("I5").Value = 0.00699771552698574
$ws.Range("J5").Value = 0.004690590957742906
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 0.68591268109
$ws.Range("R5").Value = 2.74365072436
$ws.Range("S5").Value = 0.00699771552698574
$ws.Range("T5").Value = 0.004690590957742906

# Row 6
$ws.Range("G6").Value = 272.1243743333334
$ws.Range("H6").Value = 816.3731230000001
$ws.Range("I6").Value = 0.09080481624829614
$ws.Range("J6").Value = 0.09130013538801972
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 8.900644035694668
$ws.Range("R6").Value = 53.40386421416801
$ws.Range("S6").Value = 0.09080481624829614
$ws.Range("T6").Value = 0.09130013538801972

# Row 7
$ws.Range("G7").Value = 1774.170247333333
$ws.Range("H7").Value = 5322.510741999999
$ws.Range("I7").Value = 0.5920204821673095
$ws.Range("J7").Value = 0.5952498161172183
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 58.02956044977866
$ws.Range("R7").Value = 348.177362698672
$ws.Range("S7").Value = 0.5920204821673095
$ws.Range("T7").Value = 0.5952498161172183
